$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New odds/match data fetched for rows 226-230 (id 224-228).
# Columns: B=matchId, E=date(serial), F=HomeTeam, G=AwayTeam,
# K..P=opening/closing 1X2 odds, Q=Asian handicap line,
# R..V=AH/OU odds. A, C, D, W..AA are unchanged.
$rows = @(
    @{ Row=226; B=6774877; E=45383.3125;         F="Puszcza Niepolomice"; G="Radomiak Radom";   K=2.625; L=3.4; M=2.6;  N=2.625; O=3.4;  P=2.6;  Q=0;      R=1.975; S=1.875; T=2.25; U=1.825; V=2.025 },
    @{ Row=227; B=6775576; E=45383.41666666666;  F="Stal Mielec";         G="Lech Poznan";       K=4.333; L=3.4; M=1.85; N=4.333; O=3.4;  P=1.85; Q=0.5;    R=1.975; S=1.875; T=2.25; U=1.925; V=1.925 },
    @{ Row=228; B=6775578; E=45383.52083333334;  F="Widzew Lodz";         G="Korona Kielce";     K=2.25;  L=3.2; M=3.4;  N=2.2;   O=3.25; P=3.4;  Q=-0.25;  R=1.85;  S=2;     T=2.5;  U=2.025; V=1.825 },
    @{ Row=229; B=6775573; E=45383.625;          F="Gornik Zabrze";       G="Legia Warsaw";      K=3.6;   L=3.5; M=2.05; N=3.8;   O=3.5;  P=1.95; Q=0.5;    R=1.825; S=2.025; T=2.5;  U=1.925; V=1.925 },
    @{ Row=230; B=6775577; E=45384.58333333334;  F="Warta Poznan";        G="Zaglebie Lubin";    K=2.9;   L=3.1; M=2.55; N=2.9;   O=3.1;  P=2.55; Q=0;      R=2.05;  S=1.8;   T=2.25; U=2.05;  V=1.8 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.B    # B: id
    $ws.Cells.Item($row, 5).Value = $r.E    # E: Date
    $ws.Cells.Item($row, 6).Value = $r.F    # F: HomeTeam
    $ws.Cells.Item($row, 7).Value = $r.G    # G: AwayTeam
    $ws.Cells.Item($row, 11).Value = $r.K   # K: oddH_op
    $ws.Cells.Item($row, 12).Value = $r.L   # L: oddD_op
    $ws.Cells.Item($row, 13).Value = $r.M   # M: oddA_op
    $ws.Cells.Item($row, 14).Value = $r.N   # N: oddH
    $ws.Cells.Item($row, 15).Value = $r.O   # O: oddD
    $ws.Cells.Item($row, 16).Value = $r.P   # P: oddA
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Ah
    $ws.Cells.Item($row, 18).Value = $r.R   # R: oddAHH
    $ws.Cells.Item($row, 19).Value = $r.S   # S: oddAHA
    $ws.Cells.Item($row, 20).Value = $r.T   # T: AhOU
    $ws.Cells.Item($row, 21).Value = $r.U   # U: oddAHOver
    $ws.Cells.Item($row, 22).Value = $r.V   # V: oddAHUnder
}

# The last 4 rows of stale data (old rows 231-234) are dropped entirely.
$ws.Rows("231:234").Delete()
